$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.476.61'

$ws.Range("D3").Value = '2.690.23'
$ws.Range("E3").Value = '  +1.67%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.82%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.542'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.46%  '

$ws.Range("D9").Value = '2.689.16'
$ws.Range("E9").Value = '  +1.64%  '

$ws.Range("E10").Value = '  -4.34%  '

$ws.Range("E11").Value = '  -0.81%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.29'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.358'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.54%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.26'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.56%  '

$ws.Range("D15").Value = '3.209.27'
$ws.Range("E15").Value = '  +2.55%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000186'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.68%  '

$ws.Range("D17").Value = '68.428.06'
$ws.Range("E17").Value = '  +0.27%  '

$ws.Range("D18").Value = '2.714.95'
$ws.Range("E18").Value = '  +2.31%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '365.41'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.35%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.59'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.52'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.73%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.88'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.89%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '75.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.07%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.37%  '

$ws.Range("D28").Value = '2.842.32'
$ws.Range("E28").Value = '  +1.86%  '

$ws.Range("E29").Value = '  -1.57%  '

$ws.Range("E30").Value = '  +0.10%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '578.66'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.26'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.69%  '

$ws.Range("E33").Value = '  +0.65%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.93'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.91%  '

$ws.Range("E35").Value = '  +4.85%  '

$ws.Range("E36").Value = '  +0.98%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.08%  '

$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '161.61'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.68%  '

$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.99'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.89%  '

$ws.Range("E40").Value = '  +1.75%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.90'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.43%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.39'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.68%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.86'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.48%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.63'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.31%  '

$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.02%  '

$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.0₆0316'
$ws.Range("E46").Value = '  -7.47%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '158.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.14%  '

$ws.Range("E48").Value = '  +5.01%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.76'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.80%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.600'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.87%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.04'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.03%  '
